$p = $ppt.ActivePresentation

# The deck currently carries two theme parts: the one actually driving the
# slide master / slides ("Integral" / Red Violet) and a second, unused one
# only linked from the notes master ("Office Theme" / Office). The edit
# swaps the two themes' colour schemes (font scheme and format scheme are
# already identical between them) so the deck now uses the plain "Office"
# palette.
#
# PowerPoint's Theme Colors gallery maps 1:1 onto the 12 theme colour
# slots (Dk1, Lt1, Dk2, Lt2, Accent1-6, Hyperlink, FollowedHyperlink) via
# Slide.ThemeColorScheme, so drive the swap through that object.

$officeColors = @(
    0,        # Dark 1    - 000000
    16777215, # Light 1   - FFFFFF
    6968388,  # Dark 2    - 44546A
    15132391, # Light 2   - E7E6E6
    13998939, # Accent 1  - 5B9BD5
    3243501,  # Accent 2  - ED7D31
    10855845, # Accent 3  - A5A5A5
    49407,    # Accent 4  - FFC000
    12874308, # Accent 5  - 4472C4
    4697456,  # Accent 6  - 70AD47
    12673797, # Hyperlink - 0563C1
    7491477   # Followed Hyperlink - 954F72
)

$slide = $p.Slides.Item(1)
$themeColors = $slide.ThemeColorScheme

for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Colors($i).RGB = $officeColors[$i - 1]
}
